$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-unit voltage magnitudes (vm_pu) for the 380 kV case (Case_4_207).
# Columns correspond to bus columns B,C,D,E,F,I,J,K,L,M,N (G stays 1.0, H is blank).
$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$dataText = @"
1.02,1.039071318908831,1.04398819050538,1.042724914362951,1.051837051097158,1.036623983560882,1.044165322405027,1.04676050673295,1.045500792026213,1.054587447587773,1.018504033385264
1.02,1.040478437287656,1.045068637076232,1.044081514291748,1.053137988766934,1.036957648091379,1.045215204850368,1.047651396760461,1.046666853563893,1.055699856956297,1.018863631691469
1.02,1.041387793640134,1.045766489548504,1.044958555415854,1.053978902121754,1.037171271190628,1.045892992380143,1.048226004133231,1.047420073971934,1.056418234422609,1.019095502202878
1.02,1.041769819004053,1.046059566367688,1.045327083997819,1.054332215581724,1.037260533849039,1.046177565540044,1.048467127425152,1.04773642103023,1.056719902796133,1.019192787053096
1.02,1.041833947189474,1.046108757739445,1.045388951180155,1.054391526478731,1.037275489552774,1.046225325118023,1.048507587280181,1.047789519226039,1.056770534562079,1.019209110296717
1.02,1.041392899327965,1.045770406831365,1.044963480412237,1.053983623919728,1.037172466060434,1.04589679630458,1.04822922776593,1.04742430221705,1.056422266651665,1.019096802886605
1.02,1.039547101307436,1.044353597213175,1.043183545920561,1.052276893903344,1.03673722001114,1.044520459778857,1.047061973826546,1.045895139606577,1.054963689066912,1.018625730554753
1.02,1.036285558040849,1.041847128724461,1.040040963950812,1.04926248473132,1.035952734789999,1.042083079760739,1.044990749656415,1.043190415689459,1.052382396999259,1.017789352020864
1.02,1.034104778707959,1.0401692986945,1.037941498619562,1.047247947837481,1.035417871443892,1.040449774479286,1.043600066696753,1.041380174759858,1.050653834382177,1.017227456800998
1.02,1.033158874696602,1.039441105881977,1.037031292688418,1.046374400759257,1.035183432996924,1.039740489928289,1.042995499387483,1.040594577275859,1.04990346600963,1.016983109427904
1.02,1.032807274603482,1.03917036572456,1.036693026911777,1.046049734942047,1.035095923723696,1.039476716727471,1.042770572674673,1.040302502929215,1.049624457085146,1.016892189604813
1.02,1.03288270540821,1.039228452103531,1.036765594073875,1.046119385615705,1.03511471417565,1.039533311233674,1.042818836732661,1.040365166103047,1.049684318589123,1.01691169940645
1.02,1.033129816448135,1.039418731682417,1.037003335133623,1.046347567719908,1.035176208203039,1.039718692774044,1.042976914314068,1.040570439826963,1.04988040896248,1.016975597206275
1.02,1.033282036479815,1.039535935036112,1.037149791909621,1.046488132704038,1.035214039880477,1.039832870790126,1.043074262840326,1.040696880008198,1.050001188311986,1.017014945720671
1.02,1.03416752067756,1.040217590661322,1.038001881809041,1.047305895713981,1.035433370377313,1.040496803690847,1.043640139086616,1.041432274924703,1.05070359364048,1.017243651217445
1.02,1.034722524958093,1.040644722044272,1.038536070279751,1.047818521285791,1.03557018919523,1.040912718232956,1.043994454931494,1.04189309610347,1.051143684702458,1.017386831739694
1.02,1.035046094491236,1.040893698623604,1.0388475453727,1.048117407897721,1.035649719492853,1.041155116441611,1.044200890813384,1.042161716698708,1.051400200408329,1.017470245988865
1.02,1.035156397248264,1.040978565866719,1.038953732190057,1.048219300448868,1.03567679088279,1.041237734581243,1.04427124113012,1.042253280943374,1.051487634916358,1.017498671089929
1.02,1.034662994342241,1.040598911669794,1.038478768127036,1.047763533782792,1.03555553816523,1.040868115024555,1.043956464071548,1.041843671830385,1.051096485953182,1.017371480234501
1.02,1.033057055342432,1.03936270621629,1.036933331189575,1.046280379108027,1.035158111585479,1.039664111220593,1.042930374473391,1.04050999928572,1.049822673243626,1.016956785280608
1.02,1.032045892293222,1.038583966131027,1.035960639517499,1.045346750041482,1.034905754223111,1.038905290882943,1.042283126105237,1.039669910498166,1.049020103872811,1.016695133246419
1.02,1.032582068324697,1.038996933543627,1.036476379905222,1.045841791565683,1.035039769287323,1.039307729641427,1.042626445616856,1.040115406743074,1.049445721267658,1.016833927393308
1.02,1.034689894144571,1.040619611906836,1.03850466083739,1.047788380637221,1.035562159181229,1.040888269905704,1.043973631218547,1.041866005053236,1.051117813598925,1.017378417229641
1.02,1.037129849802828,1.042496302516233,1.040854151811659,1.050042630952069,1.036157629118677,1.042714660405025,1.045527934869275,1.043890881589773,1.05305106228863,1.018006329511577
"@

$dataRows = $dataText -split "`n" | Where-Object { $_.Trim() -ne "" }

$startRow = 2
for ($i = 0; $i -lt $dataRows.Count; $i++) {
    $rowNum = $startRow + $i
    $values = $dataRows[$i] -split ","
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $colLetter = $columns[$j]
        $cellValue = [double]$values[$j]
        $ws.Range("$colLetter$rowNum").Value = $cellValue
    }
}
